# Commit: "Fixed bug in Ancas (removed old function declaration) - Removed prints"
#
# This adds three helper columns (P="a", Q="b", R="x") that interpolate
# between two bounds (a..b) using a cosine-ease curve driven by the
# existing "index" column S, for the first 16-row block (rows 2-17).
# It also drops a single leftover literal value in R23, and moves the
# active selection to AB2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings / headers -----------------------------------
# Entered in this order (x, then a, then b) so the shared-string table
# matches the authoring order of the edit.
$ws.Range("R1").Value = "x"
$ws.Range("P1").Value = "a"
$ws.Range("Q1").Value = "b"

# --- New columns P (a), Q (b), R (cosine interpolation) for rows 2-17
$ws.Range("P2:P17").Value = 0
$ws.Range("Q2:Q17").Value = 2808.9803550000001
$ws.Range("R2:R17").Formula = "=(Q2-P2)/2*COS(PI()*(15-S2+1)/15)+(Q2+P2)/2"

# --- Leftover literal value in R23 (previously produced by a formula
# that was subsequently cleaned up / removed, per the commit message)
$ws.Range("R23").Value = 2508.3142275281102

# --- Move the active selection from AE14 to AB2
$ws.Range("AB2").Select() | Out-Null
